$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.894.58'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.83%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.548.38'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.62%  '

# Row 4
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '617.76'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.71%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.98'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.56%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.543.13'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.49%  '

# Row 8
$ws.Range("E8").Value = '  -0.01%  '

# Row 9
$ws.Range("E9").Value = '  +1.61%  '

# Row 10
$ws.Range("E10").Value = '  +5.29%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.30'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +5.22%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.436'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.62%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.95'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +4.42%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.148.36'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.59%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.545.60'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.65%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.860.42'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.92%  '

# Row 18
$ws.Range("E18").Value = '  +0.25%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.75'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +5.71%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.87'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +5.69%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.94'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +10.18%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '453.29'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.52%  '

# Row 23
$ws.Range("E23").Value = '  +2.77%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.10'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.15%  '

# Row 25
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.49'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +3.43%  '

# Row 26
$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000129'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.09%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.686.27'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.47%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.11%  '

# Row 29
$ws.Range("E29").Value = '  +8.23%  '

# Row 30
$ws.Range("E30").Value = '  +3.04%  '

# Row 32
$ws.Range("E32").Value = '  +2.55%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.15%  '

# Row 34
$ws.Range("E34").Value = '  +4.37%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '26.00'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.93%  '

# Row 36
$ws.Range("E36").Value = '  +3.59%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.541.41'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.57%  '

# Row 38
$ws.Range("E38").Value = '  +3.14%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.36'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +6.84%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '178.64'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +4.23%  '

# Row 42
$ws.Range("E42").Value = '  -0.12%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0914'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +5.16%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.56'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.72%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '30.85'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +15.23%  '

# Row 46
$ws.Range("E46").Value = '  +0.95%  '

# Row 47
$ws.Range("E47").Value = '  +6.77%  '

# Row 48
$ws.Range("E48").Value = '  +0.37%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.64'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.05%  '

# Row 50
$ws.Range("E50").Value = '  +3.29%  '

# Row 51
$ws.Range("B51").Value = 'SuiNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.02'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.73%  '
